{"js": "// Apply the \"Added many more features\" edits to the Hawaiian Fever Slot\n// review document. Each entry is an exact (old, new) text pair; the title/\n// CTA line occurs twice in the document (the H1 heading and the bolded\n// line near the end), so every match returned by search() is replaced.\nconst replacements = [\n  [\n    \"Play Hawaiian Fever Slot for Free - Review & Rating 2021\",\n    \"Play Hawaiian Fever Slot for Free | Review & Bonus Features\",\n  ],\n  [\n    \"Beautiful graphics and soundtrack that transport you to Hawaii\",\n    \"Tropical Hawaiian theme with stunning graphics and soundtrack\",\n  ],\n  [\n    \"Wide range of betting options\",\n    \"Wide range of betting options for all players\",\n  ],\n  [\n    \"Expanding wilds and free spins increase potential wins\",\n    \"Exciting bonus features including Expanding Wilds and Free Spins\",\n  ],\n  [\n    \"Compatible with mobile devices\",\n    \"High maximum potential win of up to 1,860 times the bet\",\n  ],\n  [\n    \"Simplified layout may not appeal to some players\",\n    \"Gameplay may be too straightforward for some players\",\n  ],\n  [\n    \"Average RTP and high volatility may not suit all players\",\n    \"Average RTP of 95.04% may not appeal to players looking for higher returns\",\n  ],\n  [\n    \"Looking for a slot game that will transport you to Hawaii? Read our Hawaiian Fever Slot review and play for free. Expanding wilds, free spins, and more!\",\n    \"Read our review of Hawaiian Fever Slot and play for free with exciting bonus features.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Hawaiian Fever Slot\n# review document via Find/Replace. The title/CTA line occurs twice in the\n# document (the H1 heading and the bolded line near the end) \u2014 wdReplaceAll\n# (2) takes care of both, while MatchCase=$true keeps the lowercase\n# \"wide range of betting options\" phrase inside the body paragraph intact.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"Play Hawaiian Fever Slot for Free - Review & Rating 2021\", \"Play Hawaiian Fever Slot for Free | Review & Bonus Features\"),\n  @(\"Beautiful graphics and soundtrack that transport you to Hawaii\", \"Tropical Hawaiian theme with stunning graphics and soundtrack\"),\n  @(\"Wide range of betting options\", \"Wide range of betting options for all players\"),\n  @(\"Expanding wilds and free spins increase potential wins\", \"Exciting bonus features including Expanding Wilds and Free Spins\"),\n  @(\"Compatible with mobile devices\", \"High maximum potential win of up to 1,860 times the bet\"),\n  @(\"Simplified layout may not appeal to some players\", \"Gameplay may be too straightforward for some players\"),\n  @(\"Average RTP and high volatility may not suit all players\", \"Average RTP of 95.04% may not appeal to players looking for higher returns\"),\n  @(\"Looking for a slot game that will transport you to Hawaii? Read our Hawaiian Fever Slot review and play for free. Expanding wilds, free spins, and more!\", \"Read our review of Hawaiian Fever Slot and play for free with exciting bonus features.\")\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
